$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) cells ---
# Force Text number format first so strings like "1.003" or "297.62" are not
# re-interpreted as numeric values by Excel; ClearFormats afterwards restores
# the default (unstyled) cell formatting used throughout the sheet.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.294.99"
$ws.Range("D3").Value = "1.625.96"
$ws.Range("D4").Value = "1.003"
$ws.Range("D6").Value = "297.62"
$ws.Range("D8").Value = "49.94"
$ws.Range("D9").Value = "0.3460"
$ws.Range("D10").Value = "0.08012"
$ws.Range("D12").Value = "1.003"
$ws.Range("D14").Value = "6.281"
$ws.Range("D15").Value = "7.197"
$ws.Range("D16").Value = "0.00001183"
$ws.Range("D17").Value = "1.616.75"
$ws.Range("D18").Value = "94.41"
$ws.Range("D19").Value = "0.06938"
$ws.Range("D20").Value = "6.579"
$ws.Range("D21").Value = "17.20"
$ws.Range("D23").Value = "12.31"
$ws.Range("D24").Value = "23.311.31"
$ws.Range("D25").Value = "2.429"
$ws.Range("D26").Value = "2.971"
$ws.Range("D27").Value = "20.81"
$ws.Range("D28").Value = "151.38"
$ws.Range("D29").Value = "5.154"
$ws.Range("D30").Value = "131.36"
$ws.Range("D31").Value = "1.808.05"
$ws.Range("D32").Value = "6.680"
$ws.Range("D33").Value = "2.114"
$ws.Range("D34").Value = "11.15"
$ws.Range("D35").Value = "0.9700"
$ws.Range("D36").Value = "0.08741"
$ws.Range("D37").Value = "0.02645"
$ws.Range("D39").Value = "5.815"
$ws.Range("D40").Value = "0.06662"
$ws.Range("D41").Value = "12.64"
$ws.Range("D42").Value = "0.6762"
$ws.Range("D43").Value = "1.293"
$ws.Range("D44").Value = "15.37"
$ws.Range("D45").Value = "1.001"
$ws.Range("D46").Value = "0.6279"
$ws.Range("D47").Value = "2.223"
$ws.Range("D48").Value = "3.880"
$ws.Range("D49").Value = "126.37"
$ws.Range("D50").Value = "0.07616"
$ws.Range("D51").Value = "1.213"

$priceRange.ClearFormats()

# --- Update Coin name (B), Link (C) and Volume/1h change (E) cells ---
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("E9").Value = "  -4.30%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("E16").Value = "  -3.68%  "
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -4.01%  "
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("E32").Value = "  -5.52%  "
$ws.Range("E33").Value = "  -5.01%  "
$ws.Range("E34").Value = "  -7.61%  "
$ws.Range("E35").Value = "  -8.28%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E37").Value = "  -6.03%  "
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("E39").Value = "  -4.51%  "
$ws.Range("E40").Value = "  -5.03%  "
$ws.Range("E41").Value = "  -2.59%  "
$ws.Range("E42").Value = "  -3.46%  "
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("E47").Value = "  -3.90%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("E49").Value = "  -1.53%  "
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("E51").Value = "  +1.29%  "
